# Auto-generated Excel COM-interop script
# Applies numeric corrections to leve-profit tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as described by the source diff (scheduled runner data refresh).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value2 = 7936918
$ws.Range("I19").Value2 = 20408374
$ws.Range("J19").Value2 = 536.36365
$ws.Range("K19").Value2 = 20408374
$ws.Range("L19").Value2 = 536.36365
$ws.Range("M19").Value2 = -20408199
$ws.Range("N19").Value2 = -886.36365
$ws.Range("H62").Value2 = 2874.7693
$ws.Range("I62").Value2 = 2370.3333
$ws.Range("J62").Value2 = 3026.1
$ws.Range("K62").Value2 = 2370.3333
$ws.Range("L62").Value2 = 3026.1
$ws.Range("M62").Value2 = -1746.3333
$ws.Range("N62").Value2 = -4274.1
$ws.Range("H65").Value2 = 2874.7693
$ws.Range("I65").Value2 = 2370.3333
$ws.Range("J65").Value2 = 3026.1
$ws.Range("K65").Value2 = 11851.6665
$ws.Range("L65").Value2 = 15130.5
$ws.Range("M65").Value2 = -8731.666499999999
$ws.Range("N65").Value2 = -21370.5
$ws.Range("H107").Value2 = 1206.6875
$ws.Range("I107").Value2 = 1238.5385
$ws.Range("J107").Value2 = 1068.6666
$ws.Range("K107").Value2 = 1238.5385
$ws.Range("L107").Value2 = 1068.6666
$ws.Range("M107").Value2 = 681.4614999999999
$ws.Range("N107").Value2 = -4908.6666
$ws.Range("H112").Value2 = 14707650
$ws.Range("J112").Value2 = 25001304
$ws.Range("L112").Value2 = 75003912
$ws.Range("N112").Value2 = -75006128
$ws.Range("H113").Value2 = 3832.6667
$ws.Range("I113").Value2 = 4391.5
$ws.Range("J113").Value2 = 3460.111
$ws.Range("K113").Value2 = 4391.5
$ws.Range("L113").Value2 = 3460.111
$ws.Range("M113").Value2 = -1137.5
$ws.Range("N113").Value2 = -9968.111000000001
$ws.Range("H137").Value2 = 4579.909
$ws.Range("I137").Value2 = 6967.2
$ws.Range("J137").Value2 = 2590.5
$ws.Range("K137").Value2 = 20901.6
$ws.Range("L137").Value2 = 7771.5
$ws.Range("M137").Value2 = -18351.6
$ws.Range("N137").Value2 = -12871.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 35717036
$ws.Range("I2").Value2 = 57694364
$ws.Range("J2").Value2 = 3875
$ws.Range("K2").Value2 = 57694364
$ws.Range("L2").Value2 = 3875
$ws.Range("M2").Value2 = -57694251
$ws.Range("N2").Value2 = -4101
$ws.Range("H45").Value2 = 1817.6471
$ws.Range("I45").Value2 = 1306.25
$ws.Range("J45").Value2 = 10000
$ws.Range("K45").Value2 = 1306.25
$ws.Range("L45").Value2 = 10000
$ws.Range("M45").Value2 = -929.25
$ws.Range("N45").Value2 = -10754
$ws.Range("H74").Value2 = 2177.4062
$ws.Range("I74").Value2 = 1634.0416
$ws.Range("J74").Value2 = 3807.5
$ws.Range("K74").Value2 = 1634.0416
$ws.Range("L74").Value2 = 3807.5
$ws.Range("M74").Value2 = -760.0416
$ws.Range("N74").Value2 = -5555.5
$ws.Range("H77").Value2 = 2177.4062
$ws.Range("I77").Value2 = 1634.0416
$ws.Range("J77").Value2 = 3807.5
$ws.Range("K77").Value2 = 8170.208000000001
$ws.Range("L77").Value2 = 19037.5
$ws.Range("M77").Value2 = -3802.208000000001
$ws.Range("N77").Value2 = -27773.5
$ws.Range("H116").Value2 = 35717036
$ws.Range("I116").Value2 = 57694364
$ws.Range("J116").Value2 = 3875
$ws.Range("K116").Value2 = 57694364
$ws.Range("L116").Value2 = 3875
$ws.Range("M116").Value2 = -57692070
$ws.Range("N116").Value2 = -8463
$ws.Range("H122").Value2 = 3802.2222
$ws.Range("I122").Value2 = 2876.5
$ws.Range("J122").Value2 = 4542.8
$ws.Range("K122").Value2 = 8629.5
$ws.Range("L122").Value2 = 13628.4
$ws.Range("M122").Value2 = -6179.5
$ws.Range("N122").Value2 = -18528.4
$ws.Range("H132").Value2 = 2078.8108
$ws.Range("I132").Value2 = 1669.3438
$ws.Range("J132").Value2 = 4699.4
$ws.Range("K132").Value2 = 5008.0314
$ws.Range("L132").Value2 = 14098.2
$ws.Range("M132").Value2 = -2478.0314
$ws.Range("N132").Value2 = -19158.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 35717036
$ws.Range("I3").Value2 = 57694364
$ws.Range("J3").Value2 = 3875
$ws.Range("K3").Value2 = 57694364
$ws.Range("L3").Value2 = 3875
$ws.Range("M3").Value2 = -57694250
$ws.Range("N3").Value2 = -4103
$ws.Range("H132").Value2 = 30000
$ws.Range("J132").Value2 = 30000
$ws.Range("L132").Value2 = 30000
$ws.Range("N132").Value2 = -40120
$ws.Range("H140").Value2 = 35571.43
$ws.Range("J140").Value2 = 35571.43
$ws.Range("L140").Value2 = 35571.43
$ws.Range("N140").Value2 = -45931.43

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 1500
$ws.Range("I16").Value2 = 747.75
$ws.Range("J16").Value2 = 2703.6
$ws.Range("K16").Value2 = 747.75
$ws.Range("L16").Value2 = 2703.6
$ws.Range("M16").Value2 = -460.75
$ws.Range("N16").Value2 = -3277.6
$ws.Range("H23").Value2 = 82009.8
$ws.Range("I23").Value2 = 10009
$ws.Range("J23").Value2 = 100010
$ws.Range("K23").Value2 = 10009
$ws.Range("L23").Value2 = 100010
$ws.Range("M23").Value2 = -9769
$ws.Range("N23").Value2 = -100490
$ws.Range("H27").Value2 = 82009.8
$ws.Range("I27").Value2 = 10009
$ws.Range("J27").Value2 = 100010
$ws.Range("K27").Value2 = 10009
$ws.Range("L27").Value2 = 100010
$ws.Range("M27").Value2 = -9817
$ws.Range("N27").Value2 = -100394
$ws.Range("H31").Value2 = 3203.352
$ws.Range("J31").Value2 = 4338.1514
$ws.Range("L31").Value2 = 4338.1514
$ws.Range("N31").Value2 = -4928.1514
$ws.Range("H33").Value2 = 5142.857
$ws.Range("I33").Value2 = 5142.857
$ws.Range("J33").Value2 = 0
$ws.Range("K33").Value2 = 5142.857
$ws.Range("L33").Value2 = 0
$ws.Range("M33").Value2 = -4763.857
$ws.Range("N33").ClearContents()
$ws.Range("H34").Value2 = 3203.352
$ws.Range("J34").Value2 = 4338.1514
$ws.Range("L34").Value2 = 4338.1514
$ws.Range("N34").Value2 = -4742.1514
$ws.Range("H58").Value2 = 7355760
$ws.Range("I58").Value2 = 1511.262
$ws.Range("J58").Value2 = 19235700
$ws.Range("K58").Value2 = 1511.262
$ws.Range("L58").Value2 = 19235700
$ws.Range("M58").Value2 = -1308.262
$ws.Range("N58").Value2 = -19236106
$ws.Range("H94").Value2 = 5412.5
$ws.Range("I94").Value2 = 12170.667
$ws.Range("J94").Value2 = 3159.7778
$ws.Range("K94").Value2 = 12170.667
$ws.Range("L94").Value2 = 3159.7778
$ws.Range("M94").Value2 = -11719.667
$ws.Range("N94").Value2 = -4061.7778
$ws.Range("H99").Value2 = 3207.2666
$ws.Range("I99").Value2 = 1793.3334
$ws.Range("J99").Value2 = 4149.8887
$ws.Range("K99").Value2 = 1793.3334
$ws.Range("L99").Value2 = 4149.8887
$ws.Range("M99").Value2 = -295.3334
$ws.Range("N99").Value2 = -7145.8887
$ws.Range("H107").Value2 = 1004.91174
$ws.Range("I107").Value2 = 815.95654
$ws.Range("J107").Value2 = 1400
$ws.Range("K107").Value2 = 815.95654
$ws.Range("L107").Value2 = 1400
$ws.Range("M107").Value2 = 1104.04346
$ws.Range("N107").Value2 = -5240
$ws.Range("H113").Value2 = 1500
$ws.Range("I113").Value2 = 747.75
$ws.Range("J113").Value2 = 2703.6
$ws.Range("K113").Value2 = 747.75
$ws.Range("L113").Value2 = 2703.6
$ws.Range("M113").Value2 = 1422.25
$ws.Range("N113").Value2 = -7043.6
$ws.Range("H115").Value2 = 35999.5
$ws.Range("J115").Value2 = 35999.5
$ws.Range("L115").Value2 = 35999.5
$ws.Range("N115").Value2 = -38349.5
$ws.Range("H126").Value2 = 3207.2666
$ws.Range("I126").Value2 = 1793.3334
$ws.Range("J126").Value2 = 4149.8887
$ws.Range("K126").Value2 = 5380.0002
$ws.Range("L126").Value2 = 12449.6661
$ws.Range("M126").Value2 = -2910.0002
$ws.Range("N126").Value2 = -17389.6661
$ws.Range("H132").Value2 = 3404.3809
$ws.Range("I132").Value2 = 3005.875
$ws.Range("K132").Value2 = 9017.625
$ws.Range("M132").Value2 = -6487.625
$ws.Range("H136").Value2 = 7355760
$ws.Range("I136").Value2 = 1511.262
$ws.Range("J136").Value2 = 19235700
$ws.Range("K136").Value2 = 4533.786
$ws.Range("L136").Value2 = 57707100
$ws.Range("M136").Value2 = -1983.786
$ws.Range("N136").Value2 = -57712200

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value2 = 969
$ws.Range("I97").Value2 = 398.33334
$ws.Range("J97").Value2 = 1825
$ws.Range("K97").Value2 = 1195.00002
$ws.Range("L97").Value2 = 5475
$ws.Range("M97").Value2 = -699.0000199999999
$ws.Range("N97").Value2 = -6467
$ws.Range("H122").Value2 = 2186.625
$ws.Range("I122").Value2 = 790
$ws.Range("J122").Value2 = 2386.1428
$ws.Range("K122").Value2 = 7110
$ws.Range("L122").Value2 = 21475.2852
$ws.Range("M122").Value2 = -4660
$ws.Range("N122").Value2 = -26375.2852
$ws.Range("H131").Value2 = 1427.1724
$ws.Range("J131").Value2 = 1125.0197
$ws.Range("L131").Value2 = 3375.0591
$ws.Range("N131").Value2 = -13455.0591

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 2780
$ws.Range("I80").Value2 = 2567.3157
$ws.Range("J80").Value2 = 3068.6428
$ws.Range("K80").Value2 = 2567.3157
$ws.Range("L80").Value2 = 3068.6428
$ws.Range("M80").Value2 = -1569.3157
$ws.Range("N80").Value2 = -5064.6428
$ws.Range("H83").Value2 = 2780
$ws.Range("I83").Value2 = 2567.3157
$ws.Range("J83").Value2 = 3068.6428
$ws.Range("K83").Value2 = 12836.5785
$ws.Range("L83").Value2 = 15343.214
$ws.Range("M83").Value2 = -7844.5785
$ws.Range("N83").Value2 = -25327.214
$ws.Range("H102").Value2 = 62822.234
$ws.Range("I102").Value2 = 3165.75
$ws.Range("K102").Value2 = 3165.75
$ws.Range("M102").Value2 = -1543.75
$ws.Range("H122").Value2 = 5383.3403
$ws.Range("I122").Value2 = 4432.2085
$ws.Range("J122").Value2 = 6375.826
$ws.Range("K122").Value2 = 13296.6255
$ws.Range("L122").Value2 = 19127.478
$ws.Range("M122").Value2 = -10846.6255
$ws.Range("N122").Value2 = -24027.478
$ws.Range("H132").Value2 = 4472.15
$ws.Range("I132").Value2 = 4986.4346
$ws.Range("J132").Value2 = 3776.353
$ws.Range("K132").Value2 = 14959.3038
$ws.Range("L132").Value2 = 11329.059
$ws.Range("M132").Value2 = -12429.3038
$ws.Range("N132").Value2 = -16389.059

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 8478
$ws.Range("I40").Value2 = 12556
$ws.Range("J40").Value2 = 4400
$ws.Range("K40").Value2 = 12556
$ws.Range("L40").Value2 = 4400
$ws.Range("M40").Value2 = -12420
$ws.Range("N40").Value2 = -4672

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value2 = 655.9091
$ws.Range("I100").Value2 = 550
$ws.Range("J100").Value2 = 716.4286
$ws.Range("K100").Value2 = 1100
$ws.Range("L100").Value2 = 1432.8572
$ws.Range("M100").Value2 = -559
$ws.Range("N100").Value2 = -2514.8572
$ws.Range("H122").Value2 = 1980.375
$ws.Range("I122").Value2 = 1456.5
$ws.Range("K122").Value2 = 4369.5
$ws.Range("M122").Value2 = -1919.5

Write-Host "Applied $([int]271) cell updates across 8 sheets."
